# Generate Report for Handoff
# The b1e11ecd-c9be-41ac-b89a-b48322140d5d.md file has finished translation and
# is now ready for handoff: update its Status/Priority/timestamps on every sheet.

$wb = $excel.ActiveWorkbook

# The width Excel ends up storing is snapped to whole-pixel increments, so we
# pick an input that lands on the pixel bucket closest to the wider "Ready for
# handoff" / "Ready for handoff" column target.
$newStatusColWidth = 16.3333333333333

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-14 16:19:56"
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-14 16:19:48"
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-14 16:19:56"
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
